$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 2 ("H 72") and shift the remaining rows up,
# matching the observed diff where every subsequent row moves up by one
# and the used range shrinks from A1:F63 to A1:F62.
$ws.Rows.Item(2).Delete()
